$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.956.33"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.84%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.973.04"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.64%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -6.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.597"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.41%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "53.57"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -7.14%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "59.84"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.67%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.367"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.52%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0747"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -7.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0978"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.31%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.264.54"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "13.84"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.97%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.71"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.66%  "
$ws.Range("E16").Value = "  -9.65%  "
$ws.Range("E17").Value = "  -7.75%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.966.33"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "36.848.53"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.88%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.98"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.10%  "
$ws.Range("E21").Value = "  -6.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "226.93"
$ws.Range("D22").Style = "Normal"
$ws.Range("E23").Value = "  -6.62%  "
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.34"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.69%  "
$ws.Range("E26").Value = "  -12.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "160.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.52%  "
$ws.Range("E28").Value = "  -6.94%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.95"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.23%  "
$ws.Range("E30").Value = "  -11.03%  "
$ws.Range("E31").Value = "  -5.66%  "
$ws.Range("E32").Value = "  -3.89%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.36"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -8.99%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0611"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -8.56%  "
$ws.Range("E35").Value = "  -7.30%  "
$ws.Range("E36").Value = "  -8.07%  "
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.79"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.29"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.15"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.62%  "
$ws.Range("E41").Value = "  -0.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.422.93"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.90%  "
$ws.Range("E43").Value = "  -6.81%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0887"
$ws.Range("D44").Style = "Normal"
$ws.Range("E45").Value = "  -7.64%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "86.83"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.45%  "
$ws.Range("E47").Value = "  -7.59%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.87"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +18.63%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.987"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.55%  "
$ws.Range("E50").Value = "  -1.12%  "
$ws.Range("E51").Value = "  -11.67%  "
